$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list (column D, rows 25-30)
$ws.Range("D25").Value = 1252.84
$ws.Range("D26").Value = 1396.94
$ws.Range("D27").Value = 1574.396
$ws.Range("D28").Value = 1723.831
$ws.Range("D29").Value = 1894.6
$ws.Range("D30").Value = 2068.056
